# chore: update Sheets via scheduled runner
#
# Refreshes Hades_Profits.xlsx leve-profit rows with the latest market-board
# average prices (currentAveragePrice / NQ / HQ columns H/I/J) and the
# resulting LevePriceNQ/HQ (K/L) + LeveProfitNQ/HQ (M/N) recomputations,
# one worksheet per crafting job.
$wb = $excel.ActiveWorkbook

# ALC ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 133 - Big Brush, Big Dreams / Ginseng Angle Brush
$ws.Range("H133").Value = 29940
$ws.Range("J133").Value = 29940
$ws.Range("L133").Value = 29940
$ws.Range("N133").Value = -40060

# Row 140 - Tome for Tradition / Book of Ra'Kaznar
$ws.Range("H140").Value = 55277.5
$ws.Range("J140").Value = 55277.5
$ws.Range("L140").Value = 55277.5
$ws.Range("N140").Value = -65637.5

# ARM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 45 - Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 1674.5
$ws.Range("I45").Value = 2133.5715
$ws.Range("K45").Value = 2133.5715
$ws.Range("M45").Value = -1756.5715

# Row 132 - Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 52321.414
$ws.Range("I132").Value = 31111.092
$ws.Range("K132").Value = 93333.276
$ws.Range("M132").Value = -90803.276

# Row 134 - Brace for More Vambraces / Ruthenium Vambraces of Maiming
$ws.Range("H134").Value = 46913.57
$ws.Range("J134").Value = 46913.57
$ws.Range("L134").Value = 46913.57
$ws.Range("N134").Value = -57053.57

# BSM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 139 - Maul Me / Titanium Gold Maul
$ws.Range("H139").Value = 50720.832
$ws.Range("J139").Value = 40780
$ws.Range("L139").Value = 40780
$ws.Range("N139").Value = -51060

# Row 140 - Ceremonial Teeth / Ra'Kaznar Twinfangs
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

# CRP ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 31 - Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 5788.1816
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 5788.1816
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 5788.1816
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -6378.1816

# Row 34 - Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 5788.1816
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 5788.1816
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 5788.1816
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -6192.1816

# Row 100 - Run Before They Walk / Pine Cane
$ws.Range("H100").Value = 37000
$ws.Range("J100").Value = 37000
$ws.Range("L100").Value = 37000
$ws.Range("N100").Value = -39164

# Row 135 - The Wing's Wings / Ceiba Wings
$ws.Range("H135").Value = 49800
$ws.Range("J135").Value = 49800
$ws.Range("L135").Value = 49800
$ws.Range("N135").Value = -59940

# Row 139 - Weaving a Path / Acacia Spinning Wheel
$ws.Range("H139").Value = 58694
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 58694
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 58694
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -68974

# Row 140 - Spear Pressure / Claro Walnut Spear
$ws.Range("H140").Value = 51000
$ws.Range("J140").Value = 51000
$ws.Range("L140").Value = 51000
$ws.Range("N140").Value = -61360

# CUL ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 22 - A Total Nut Job / Walnut Bread
$ws.Range("H22").Value = 1216.3334
$ws.Range("J22").Value = 1600
$ws.Range("L22").Value = 4800
$ws.Range("N22").Value = -5138

# Row 27 - Brain Food / Walnut Bread
$ws.Range("H27").Value = 1216.3334
$ws.Range("J27").Value = 1600
$ws.Range("L27").Value = 4800
$ws.Range("N27").Value = -5004

# Row 33 - Cooking with Gas / Chicken Stock
$ws.Range("H33").Value = 241.44444
$ws.Range("I33").Value = 232
$ws.Range("J33").Value = 256.2857
$ws.Range("K33").Value = 1392
$ws.Range("L33").Value = 1537.7142
$ws.Range("M33").Value = -1109
$ws.Range("N33").Value = -2103.7142

# Row 39 - Bloody Good Tart, This / Blood Currant Tart
$ws.Range("H39").Value = 1625
$ws.Range("J39").Value = 1625
$ws.Range("L39").Value = 4875
$ws.Range("N39").Value = -5463

# Row 49 - Leek Soup for the Soul / Cawl Cennin
$ws.Range("H49").Value = 2700

# Row 76 - Old Victories, New Tastes / Dhalmel Fricassee
$ws.Range("H76").Value = 3200
$ws.Range("I76").Value = 1400
$ws.Range("J76").Value = 3800
$ws.Range("K76").Value = 4200
$ws.Range("L76").Value = 11400
$ws.Range("M76").Value = -3817
$ws.Range("N76").Value = -12166

# Row 79 - The Eats of Authenticity (L) / Dhalmel Fricassee
$ws.Range("H79").Value = 3200
$ws.Range("I79").Value = 1400
$ws.Range("J79").Value = 3800
$ws.Range("K79").Value = 4200
$ws.Range("L79").Value = 11400
$ws.Range("M79").Value = -2874
$ws.Range("N79").Value = -14052

# Row 106 - Herky Jerky / Jerked Jhammel
$ws.Range("H106").Value = 3500.25
$ws.Range("J106").Value = 3500.25
$ws.Range("L106").Value = 10500.75
$ws.Range("N106").Value = -12392.75

# Row 109 - Cure for What Ails / Purple Carrot Juice
$ws.Range("H109").Value = 3059.279
$ws.Range("I109").Value = 1844.9
$ws.Range("J109").Value = 3427.2727
$ws.Range("K109").Value = 5534.700000000001
$ws.Range("L109").Value = 10281.8181
$ws.Range("M109").Value = -4494.700000000001
$ws.Range("N109").Value = -12361.8181

# Row 112 - Sweet Tooth / Caramels
$ws.Range("H112").Value = 2969.842
$ws.Range("I112").Value = 2078.375
$ws.Range("J112").Value = 3618.182
$ws.Range("K112").Value = 6235.125
$ws.Range("L112").Value = 10854.546
$ws.Range("M112").Value = -5127.125
$ws.Range("N112").Value = -13070.546

# Row 113 - Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 561.15094
$ws.Range("I113").Value = 492.36365
$ws.Range("J113").Value = 609.9677
$ws.Range("K113").Value = 1477.09095
$ws.Range("L113").Value = 1829.9031
$ws.Range("M113").Value = 692.90905
$ws.Range("N113").Value = -6169.9031

# Row 131 - The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 972.0465
$ws.Range("J131").Value = 1032.974
$ws.Range("L131").Value = 3098.922
$ws.Range("N131").Value = -13178.922

# GSM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 98 - Cutting Deals / Durium Smallsword
$ws.Range("H98").Value = 29900
$ws.Range("J98").Value = 29900
$ws.Range("L98").Value = 29900
$ws.Range("N98").Value = -35890

# LTW ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 36 - Campaign in the Membrane / Toadskin Jacket
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

# Row 98 - Try Tricorne Again / Tigerskin Tricorne of Aiming
$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990

# Row 132 - Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 46831.44
$ws.Range("I132").Value = 26246.285
$ws.Range("J132").Value = 80084.38
$ws.Range("K132").Value = 78738.855
$ws.Range("L132").Value = 240253.14
$ws.Range("M132").Value = -76208.855
$ws.Range("N132").Value = -245313.14

# Row 141 - Just Generally Freezing / Gargantuaskin Trousers of Striking
$ws.Range("H141").Value = 38571.332
$ws.Range("J141").Value = 38571.332
$ws.Range("L141").Value = 38571.332
$ws.Range("N141").Value = -48931.332

# WVR ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 136 - Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 55326.594
$ws.Range("I136").Value = 31390.424
$ws.Range("K136").Value = 94171.272
$ws.Range("M136").Value = -91621.272

# Row 137 - Traditional Trousers / Sarcenet Slops of Aiming
$ws.Range("H137").Value = 45966.668
$ws.Range("J137").Value = 45966.668
$ws.Range("L137").Value = 45966.668
$ws.Range("N137").Value = -56166.668

# Row 140 - Glamorous Gloves / Thunderyards Silk Gloves of Casting
$ws.Range("H140").Value = 50397.5
$ws.Range("J140").Value = 50397.5
$ws.Range("L140").Value = 50397.5
$ws.Range("N140").Value = -60757.5
